# Apply cryptocurrency price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are plain text in the source data (e.g. "26.008.00" or
# "214.50"). Excel auto-coerces a Value assignment that *looks* numeric into a
# real number (dropping significant trailing zeros), so force text formatting for
# the duration of the write, then drop the formatting override again so the cell
# is left exactly as it was styled before (general, unstyled) but keeps its text.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.008.00"
$ws.Range("E2").Value = "  -0.36%  "
Set-TextValue "D3" "1.632.40"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "214.50"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("E9").Value = "  -3.44%  "
Set-TextValue "D10" "18.23"
$ws.Range("E10").Value = "  -7.38%  "
Set-TextValue "D11" "0.0790"
$ws.Range("E11").Value = "  -0.90%  "
Set-TextValue "D12" "1.859.43"
$ws.Range("E12").Value = "  -0.93%  "
Set-TextValue "D13" "1.658.26"
$ws.Range("E13").Value = "  -0.18%  "
Set-TextValue "D14" "4.17"
$ws.Range("E14").Value = "  -2.93%  "
Set-TextValue "D15" "0.525"
$ws.Range("E15").Value = "  -3.89%  "
Set-TextValue "D16" "26.002.46"
$ws.Range("E16").Value = "  -0.85%  "
Set-TextValue "D17" "0.0₃0740"
Set-TextValue "D18" "61.32"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("E19").Value = "  +0.17%  "
Set-TextValue "D20" "190.68"
$ws.Range("E20").Value = "  -3.12%  "
Set-TextValue "D21" "4.24"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").Value = "  +0.02%  "
Set-TextValue "D25" "144.25"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  -0.06%  "
Set-TextValue "D28" "6.77"
$ws.Range("E28").Value = "  -2.28%  "
Set-TextValue "D29" "15.18"
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("E30").Value = "  -1.81%  "
Set-TextValue "D31" "0.0480"
$ws.Range("E31").Value = "  -3.63%  "
Set-TextValue "D32" "3.13"
$ws.Range("E32").Value = "  -4.51%  "
Set-TextValue "D33" "3.12"
$ws.Range("E33").Value = "  -5.48%  "
Set-TextValue "D34" "2.41"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  -3.95%  "
Set-TextValue "D36" "1.128.67"
$ws.Range("E36").Value = "  -0.72%  "
Set-TextValue "D37" "0.859"
$ws.Range("E37").Value = "  -5.39%  "
$ws.Range("E38").Value = "  -1.16%  "
Set-TextValue "D39" "0.517"
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("E40").Value = "  -2.01%  "
Set-TextValue "D41" "98.35"
$ws.Range("E41").Value = "  -1.23%  "
Set-TextValue "D42" "0.774"
$ws.Range("E42").Value = "  -3.15%  "
Set-TextValue "D43" "1.769.71"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("E44").Value = "  -5.30%  "
Set-TextValue "D45" "0.0₆0116"
$ws.Range("E45").Value = "  -1.23%  "
Set-TextValue "D46" "54.70"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("E47").Value = "  -0.03%  "
Set-TextValue "D48" "1.49"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  +0.12%  "
Set-TextValue "D51" "7.48"
$ws.Range("E51").Value = "  -3.55%  "
